$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 649.9355
$ws.Range("I80").Value = 429.375
$ws.Range("J80").Value = 885.2
$ws.Range("K80").Value = 1288.125
$ws.Range("L80").Value = 2655.6
$ws.Range("M80").Value = -290.125
$ws.Range("N80").Value = -4651.6

$ws.Range("H83").Value = 649.9355
$ws.Range("I83").Value = 429.375
$ws.Range("J83").Value = 885.2
$ws.Range("K83").Value = 3864.375
$ws.Range("L83").Value = 7966.8
$ws.Range("M83").Value = 1127.625
$ws.Range("N83").Value = -17950.8

$ws.Range("H129").Value = 1039.8182
$ws.Range("J129").Value = 1107.6938
$ws.Range("L129").Value = 3323.0814
$ws.Range("N129").Value = -13323.0814

$ws.Range("H137").Value = 1520.2972
$ws.Range("J137").Value = 1518.5
$ws.Range("L137").Value = 4555.5
$ws.Range("N137").Value = -9655.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4614.6406
$ws.Range("I32").Value = 3091.151
$ws.Range("J32").Value = 11955.091
$ws.Range("K32").Value = 3091.151
$ws.Range("L32").Value = 11955.091
$ws.Range("M32").Value = -2804.151
$ws.Range("N32").Value = -12529.091

$ws.Range("H43").Value = 21859
$ws.Range("I43").Value = 15200
$ws.Range("J43").Value = 25188.5
$ws.Range("K43").Value = 15200
$ws.Range("L43").Value = 25188.5
$ws.Range("M43").Value = -14887
$ws.Range("N43").Value = -25814.5

$ws.Range("H45").Value = 1890.0555
$ws.Range("I45").Value = 1188.8125
$ws.Range("K45").Value = 1188.8125
$ws.Range("M45").Value = -811.8125

$ws.Range("H74").Value = 1264.6666
$ws.Range("I74").Value = 1297.6316
$ws.Range("J74").Value = 1139.4
$ws.Range("K74").Value = 1297.6316
$ws.Range("L74").Value = 1139.4
$ws.Range("M74").Value = -423.6315999999999
$ws.Range("N74").Value = -2887.4

$ws.Range("H77").Value = 1264.6666
$ws.Range("I77").Value = 1297.6316
$ws.Range("J77").Value = 1139.4
$ws.Range("K77").Value = 6488.157999999999
$ws.Range("L77").Value = 5697
$ws.Range("M77").Value = -2120.157999999999
$ws.Range("N77").Value = -14433

$ws.Range("H122").Value = 15875569
$ws.Range("I122").Value = 30304740
$ws.Range("K122").Value = 90914220
$ws.Range("M122").Value = -90911770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4679.1406
$ws.Range("I31").Value = 952.14813
$ws.Range("J31").Value = 6966.159
$ws.Range("K31").Value = 952.14813
$ws.Range("L31").Value = 6966.159
$ws.Range("M31").Value = -657.14813
$ws.Range("N31").Value = -7556.159

$ws.Range("H34").Value = 4679.1406
$ws.Range("I34").Value = 952.14813
$ws.Range("J34").Value = 6966.159
$ws.Range("K34").Value = 952.14813
$ws.Range("L34").Value = 6966.159
$ws.Range("M34").Value = -750.14813
$ws.Range("N34").Value = -7370.159

$ws.Range("H132").Value = 31562.03
$ws.Range("I132").Value = 36279.207
$ws.Range("J132").Value = 4202.4
$ws.Range("K132").Value = 108837.621
$ws.Range("L132").Value = 12607.2
$ws.Range("M132").Value = -106307.621
$ws.Range("N132").Value = -17667.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2205.6
$ws.Range("I3").Value = 998.75
$ws.Range("J3").Value = 7033
$ws.Range("K3").Value = 2996.25
$ws.Range("L3").Value = 21099
$ws.Range("M3").Value = -2884.25
$ws.Range("N3").Value = -21323

$ws.Range("H140").Value = 36202.516
$ws.Range("I140").Value = 56560.723
$ws.Range("J140").Value = 2889.0908
$ws.Range("K140").Value = 169682.169
$ws.Range("L140").Value = 8667.2724
$ws.Range("M140").Value = -164502.169
$ws.Range("N140").Value = -19027.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13565.1
$ws.Range("I70").Value = 17693.285
$ws.Range("J70").Value = 3932.6667
$ws.Range("K70").Value = 17693.285
$ws.Range("L70").Value = 3932.6667
$ws.Range("M70").Value = -17423.285
$ws.Range("N70").Value = -4472.6667

$ws.Range("H73").Value = 13565.1
$ws.Range("I73").Value = 17693.285
$ws.Range("J73").Value = 3932.6667
$ws.Range("K73").Value = 17693.285
$ws.Range("L73").Value = 3932.6667
$ws.Range("M73").Value = -16757.285
$ws.Range("N73").Value = -5804.6667

$ws.Range("H113").Value = 1241.1
$ws.Range("I113").Value = 1263.875
$ws.Range("J113").Value = 1150
$ws.Range("K113").Value = 1263.875
$ws.Range("L113").Value = 1150
$ws.Range("M113").Value = 906.125
$ws.Range("N113").Value = -5490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4230.096
$ws.Range("I132").Value = 5305.273
$ws.Range("K132").Value = 15915.819
$ws.Range("M132").Value = -13385.819

$ws.Range("H136").Value = 20838100
$ws.Range("I136").Value = 5013
$ws.Range("J136").Value = 41671188
$ws.Range("K136").Value = 15039
$ws.Range("L136").Value = 125013564
$ws.Range("M136").Value = -12489
$ws.Range("N136").Value = -125018664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 55000
$ws.Range("J60").Value = 55000
$ws.Range("L60").Value = 55000
$ws.Range("N60").Value = -56644

$ws.Range("H122").Value = 31197.205
$ws.Range("I122").Value = 38414.668
$ws.Range("K122").Value = 115244.004
$ws.Range("M122").Value = -112794.004

$ws.Range("H132").Value = 1847.2941
$ws.Range("I132").Value = 1421.6428
$ws.Range("J132").Value = 2365.4783
$ws.Range("K132").Value = 4264.928400000001
$ws.Range("L132").Value = 7096.4349
$ws.Range("M132").Value = -1734.928400000001
$ws.Range("N132").Value = -12156.4349

$ws.Range("H135").Value = 76979
$ws.Range("J135").Value = 76979
$ws.Range("L135").Value = 76979
$ws.Range("N135").Value = -87119

$ws.Range("H136").Value = 6948777.5
$ws.Range("I136").Value = 19608652
$ws.Range("J136").Value = 6265.3228
$ws.Range("K136").Value = 58825956
$ws.Range("L136").Value = 18795.9684
$ws.Range("M136").Value = -58823406
$ws.Range("N136").Value = -23895.9684

$ws.Range("H137").Value = 66253
$ws.Range("J137").Value = 66253
$ws.Range("L137").Value = 66253
$ws.Range("N137").Value = -76453

$ws.Range("H140").Value = 52980
$ws.Range("J140").Value = 52980
$ws.Range("L140").Value = 52980
$ws.Range("N140").Value = -63340
